# Apply the commit's edits to the document.
#
# Original paragraphs:
#   1. "Test me Konrad I am naught docx file"
#   2. (empty)
#   3. "Second line"
#   4. (empty)
#   5. "Third line"
#   6. (empty)
#   7. (empty)
#
# Target paragraphs:
#   1. "Lets see if sth is differen"
#   2. "or noot"                         (new paragraph)
#   3. (empty)
#   4. "or yes"
#   5. (empty)
#   6. (empty)
#
# Text is swapped in by inserting a brand-new paragraph with the new
# wording and deleting the old one (rather than overwriting the run's
# text in place) so each surviving run keeps its original/expected
# formatting shell.

$d = $word.ActiveDocument

# --- 1. "Test me Konrad I am naught docx file" -> "Lets see if sth is differen"
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphBefore()
$d.Paragraphs(1).Range.InsertBefore("Lets see if sth is differen")
$d.Paragraphs(2).Range.Delete()

# --- 2. Insert new paragraph "or noot" right after paragraph 1
$d.Paragraphs(1).Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.InsertBefore("or noot")

# --- 3. "Second line" -> "or yes"
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphBefore()
$d.Paragraphs(4).Range.InsertBefore("or yes")
$d.Paragraphs(5).Range.Delete()

# --- 4. Remove the blank paragraph and the "Third line" paragraph that used
#        to trail "Second line" (they now trail "or yes" at index 5).
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(5).Range.Delete()
